# Update "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets.
# These two sheets list the same set of events (全部类型 has one extra leading
# row compared to 展览), so the F-column (想去人数) numbers are bumped by the
# same deltas on both sheets.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row => new value, for sheet "展览"
$exhibitUpdates = @{
    8  = 627
    9  = 1073
    12 = 5010
    15 = 196
    16 = 6
    18 = 4277
    19 = 197
    20 = 1143
    21 = 116
    23 = 210
    24 = 51
    25 = 155
    26 = 56
    27 = 146
    31 = 66
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Row => new value, for sheet "全部类型" (rows shifted by +1 vs 展览)
$allUpdates = @{
    9  = 627
    10 = 1073
    13 = 5010
    16 = 196
    17 = 6
    19 = 4277
    20 = 197
    21 = 1143
    22 = 116
    24 = 210
    25 = 51
    26 = 155
    27 = 56
    28 = 146
    32 = 66
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
